# Scheduled-runner refresh of market-price-derived columns (H-N) across
# the Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Recomputed currentAveragePrice / LevePrice / LeveProfit figures; a few
# rows gain or lose their HQ-profit (N) or NQ-profit (M) cell entirely
# when that side of the recipe has no HQ/NQ variant.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 654358.6
$ws.Range("J28").Value = 894.63635
$ws.Range("L28").Value = 894.63635
$ws.Range("N28").Value = -1864.63635

$ws.Range("H98").Value = 448735.4
$ws.Range("I98").Value = 509547.03
$ws.Range("J98").Value = 2783.3333
$ws.Range("K98").Value = 509547.03
$ws.Range("L98").Value = 2783.3333
$ws.Range("M98").Value = -508049.03
$ws.Range("N98").Value = -5779.3333

$ws.Range("H122").Value = 448735.4
$ws.Range("I122").Value = 509547.03
$ws.Range("J122").Value = 2783.3333
$ws.Range("K122").Value = 1528641.09
$ws.Range("L122").Value = 8349.999899999999
$ws.Range("M122").Value = -1526191.09
$ws.Range("N122").Value = -13249.9999

$ws.Range("H132").Value = 169427.94
$ws.Range("I132").Value = 206584.94
$ws.Range("J132").Value = 32411.5
$ws.Range("K132").Value = 619754.8200000001
$ws.Range("L132").Value = 97234.5
$ws.Range("M132").Value = -617224.8200000001
$ws.Range("N132").Value = -102294.5

$ws.Range("H138").Value = 6914396
$ws.Range("I138").Value = 1467018.1
$ws.Range("J138").Value = 9093347
$ws.Range("K138").Value = 4401054.300000001
$ws.Range("L138").Value = 27280041
$ws.Range("M138").Value = -4395914.300000001
$ws.Range("N138").Value = -27290321

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H92").Value = 30000
$ws.Range("J92").Value = 30000
$ws.Range("L92").Value = 30000
$ws.Range("N92").Value = -34992

$ws.Range("H122").Value = 1704.6364
$ws.Range("I122").Value = 1616.5555
$ws.Range("J122").Value = 2101
$ws.Range("K122").Value = 4849.666499999999
$ws.Range("L122").Value = 6303
$ws.Range("M122").Value = -2399.666499999999
$ws.Range("N122").Value = -11203

$ws.Range("H139").Value = 54666.668
$ws.Range("J139").Value = 54666.668
$ws.Range("L139").Value = 54666.668
$ws.Range("N139").Value = -64946.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H108").Value = 30000
$ws.Range("J108").Value = 30000
$ws.Range("L108").Value = 30000
$ws.Range("N108").Value = -37680

$ws.Range("H134").Value = 3898.375
$ws.Range("I134").Value = 2925.2
$ws.Range("J134").Value = 5520.3335
$ws.Range("K134").Value = 8775.599999999999
$ws.Range("L134").Value = 16561.0005
$ws.Range("M134").Value = -6240.599999999999
$ws.Range("N134").Value = -21631.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 118.833336
$ws.Range("I22").Value = 128.25
$ws.Range("J22").Value = 100
$ws.Range("K22").Value = 128.25
$ws.Range("L22").Value = 100
$ws.Range("M22").Value = 221.75
$ws.Range("N22").Value = -800

$ws.Range("H31").Value = 1619.3793
$ws.Range("I31").Value = 1193.92
$ws.Range("J31").Value = 4278.5
$ws.Range("K31").Value = 1193.92
$ws.Range("L31").Value = 4278.5
$ws.Range("M31").Value = -898.9200000000001
$ws.Range("N31").Value = -4868.5

$ws.Range("H34").Value = 1619.3793
$ws.Range("I34").Value = 1193.92
$ws.Range("J34").Value = 4278.5
$ws.Range("K34").Value = 1193.92
$ws.Range("L34").Value = 4278.5
$ws.Range("M34").Value = -991.9200000000001
$ws.Range("N34").Value = -4682.5

$ws.Range("H58").Value = 1840.762
$ws.Range("I58").Value = 780.4286
$ws.Range("K58").Value = 780.4286
$ws.Range("M58").Value = -577.4286

$ws.Range("H132").Value = 2006.871
$ws.Range("I132").Value = 1238.0435
$ws.Range("J132").Value = 4217.25
$ws.Range("K132").Value = 3714.1305
$ws.Range("L132").Value = 12651.75
$ws.Range("M132").Value = -1184.1305
$ws.Range("N132").Value = -17711.75

$ws.Range("H136").Value = 1840.762
$ws.Range("I136").Value = 780.4286
$ws.Range("K136").Value = 2341.2858
$ws.Range("M136").Value = 208.7142000000003

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1028.125
$ws.Range("J5").Value = 2338.5
$ws.Range("L5").Value = 7015.5
$ws.Range("N5").Value = -7239.5

$ws.Range("H122").Value = 714.5714
$ws.Range("I122").Value = 440.6
$ws.Range("J122").Value = 1399.5
$ws.Range("K122").Value = 3965.4
$ws.Range("L122").Value = 12595.5
$ws.Range("M122").Value = -1515.4
$ws.Range("N122").Value = -17495.5

$ws.Range("H132").Value = 1300.4375
$ws.Range("I132").Value = 1152
$ws.Range("J132").Value = 1389.5
$ws.Range("K132").Value = 10368
$ws.Range("L132").Value = 12505.5
$ws.Range("M132").Value = -7838
$ws.Range("N132").Value = -17565.5

$ws.Range("H135").Value = 1028.125
$ws.Range("J135").Value = 2338.5
$ws.Range("L135").Value = 21046.5
$ws.Range("N135").Value = -26116.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 31333.334
$ws.Range("J62").Value = 31333.334
$ws.Range("L62").Value = 31333.334
$ws.Range("N62").Value = -32705.334

$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

$ws.Range("H65").Value = 31333.334
$ws.Range("J65").Value = 31333.334
$ws.Range("L65").Value = 94000.00199999999
$ws.Range("N65").Value = -100864.002

$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

$ws.Range("H107").Value = 266.46155
$ws.Range("I107").Value = 229.28572
$ws.Range("K107").Value = 229.28572
$ws.Range("M107").Value = 1690.71428

$ws.Range("H113").Value = 1366.3334
$ws.Range("I113").Value = 1339.6
$ws.Range("K113").Value = 1339.6
$ws.Range("M113").Value = 830.4000000000001

$ws.Range("H122").Value = 1236413
$ws.Range("I122").Value = 1588988.2
$ws.Range("J122").Value = 2399.5
$ws.Range("K122").Value = 4766964.6
$ws.Range("L122").Value = 7198.5
$ws.Range("M122").Value = -4764514.6
$ws.Range("N122").Value = -12098.5

$ws.Range("H132").Value = 4004.5715
$ws.Range("I132").Value = 2478
$ws.Range("K132").Value = 7434
$ws.Range("M132").Value = -4904

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 50000
$ws.Range("J2").Value = 50000
$ws.Range("L2").Value = 50000
$ws.Range("N2").Value = -50224

$ws.Range("H68").Value = 1613
$ws.Range("I68").Value = 1660.8
$ws.Range("J68").Value = 1533.3334
$ws.Range("K68").Value = 1660.8
$ws.Range("L68").Value = 1533.3334
$ws.Range("M68").Value = -911.8
$ws.Range("N68").Value = -3031.3334

$ws.Range("H71").Value = 1613
$ws.Range("I71").Value = 1660.8
$ws.Range("J71").Value = 1533.3334
$ws.Range("K71").Value = 8304
$ws.Range("L71").Value = 7666.666999999999
$ws.Range("M71").Value = -4560
$ws.Range("N71").Value = -15154.667

$ws.Range("H82").Value = 73599.86
$ws.Range("I82").Value = 101879.8
$ws.Range("J82").Value = 2900
$ws.Range("K82").Value = 101879.8
$ws.Range("L82").Value = 2900
$ws.Range("M82").Value = -101518.8
$ws.Range("N82").Value = -3622

$ws.Range("H85").Value = 73599.86
$ws.Range("I85").Value = 101879.8
$ws.Range("J85").Value = 2900
$ws.Range("K85").Value = 101879.8
$ws.Range("L85").Value = 2900
$ws.Range("M85").Value = -100631.8
$ws.Range("N85").Value = -5396

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 36375944
$ws.Range("I2").Value = 66675332
$ws.Range("J2").Value = 16680.6
$ws.Range("K2").Value = 66675332
$ws.Range("L2").Value = 16680.6
$ws.Range("M2").Value = -66675220
$ws.Range("N2").Value = -16904.6

$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()

$ws.Range("H136").Value = 12384160
$ws.Range("I136").Value = 14536982
$ws.Range("J136").Value = 5435.5
$ws.Range("K136").Value = 43610946
$ws.Range("L136").Value = 16306.5
$ws.Range("M136").Value = -43608396
$ws.Range("N136").Value = -21406.5
